# Update: Add data for 2021-11-03 (carjacking arrests by month, YoY)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet / update its "through date" label (10-25 -> 10-26)
$ws.Name = "Through 2021-10-26"

# Row 12: "October (through 10-25/26)" data row
$ws.Range("A12").Value = "October (through 10-26)"
$ws.Range("C12").Value = 23
$ws.Range("D12").Value = 0.08
$ws.Range("E12").Value = 5
$ws.Range("F12").Value = 39
$ws.Range("G12").Value = 0.1136
$ws.Range("I12").Value = 51
$ws.Range("J12").Value = 0.1356
$ws.Range("L12").Value = 51
$ws.Range("M12").Value = 0.0893
$ws.Range("O12").Value = 43
$ws.Range("P12").Value = 0.0851
$ws.Range("R12").Value = 126
$ws.Range("U12").Value = 165

# Row 13: "Total" row
$ws.Range("C13").Value = 219
$ws.Range("D13").Value = 0.1275
$ws.Range("E13").Value = 51
$ws.Range("F13").Value = 422
$ws.Range("G13").Value = 0.1078
$ws.Range("I13").Value = 628
$ws.Range("J13").Value = 0.0845
$ws.Range("L13").Value = 538
$ws.Range("M13").Value = 0.1093
$ws.Range("O13").Value = 422
$ws.Range("P13").Value = 0.1002
$ws.Range("R13").Value = 974
$ws.Range("S13").Value = 0.0516
$ws.Range("U13").Value = 1330
$ws.Range("V13").Value = 0.0581
